$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("report")
$ws2 = $wb.Worksheets.Item("service_tables")
$ws2.Columns.Item(5).Copy()
$ws1.Columns.Item(4).PasteSpecial(8)
